# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Inserts the new worker "ANDERSON DE JESUS CABRALES PADILLA" (1047436296)
# right after the already-present rows (as row 19), and re-sorts
# ANGEL MANUEL ROMERO COTA's overdue periods into ascending order
# (2201..2208) while updating his Salario Basico (column G) to 877803.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: new entry for ANDERSON DE JESUS CABRALES PADILLA
$ws.Range("C19").Value = "1047436296"
$ws.Range("D19").Value = "ANDERSON DE JESUS CABRALES PADILLA"
$ws.Range("E19").Value = "2006"
$ws.Range("F19").Value = 35112
$ws.Range("G19").Value = 877803

# Rows 20-27: ANGEL MANUEL ROMERO COTA, periods now ascending 2201..2208
$periods = @("2201","2202","2203","2204","2205","2206","2207","2208")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 20 + $i
    $ws.Range("C$row").Value = "73208312"
    $ws.Range("D$row").Value = "ANGEL MANUEL ROMERO COTA"
    $ws.Range("E$row").Value = $periods[$i]
    if ($row -eq 27) {
        $ws.Range("F$row").Value = 25749
    } else {
        $ws.Range("F$row").Value = 36341
    }
    $ws.Range("G$row").Value = 877803
}
